$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 17) with the latest analysis snapshot,
# matching the existing header layout (timestamp + 26 percentage columns a-z).
$row = 17

$ws.Range("A$row").Value = "2024-09-04 16:57:08"

$values = @{
    "B" = 0
    "C" = 0
    "D" = 0
    "E" = 0
    "F" = 33.33333333333333
    "G" = 0
    "H" = 0
    "I" = 0
    "J" = 0
    "K" = 0
    "L" = 0
    "M" = 50
    "N" = 0
    "O" = 50
    "P" = 0
    "Q" = 100
    "R" = 0
    "S" = 0
    "T" = 0
    "U" = 0
    "V" = 0
    "W" = 0
    "X" = 100
    "Y" = 0
    "Z" = 0
    "AA" = 100
}

foreach ($col in $values.Keys) {
    $ws.Range("$col$row").Value = $values[$col]
}
